$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.693.06"
$ws.Range("D3").Value = "2.036.95"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.51"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0835"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "2.337.14"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("E15").Value = "  +4.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.770"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "2.038.18"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "37.681.75"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").Value = "0.0₃0821"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.20%  "
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("E32").Value = "  +8.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("E37").Value = "  +3.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.67%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.19%  "
$ws.Range("D41").Value = "1.526.97"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0215"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0906"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("E50").Value = "  -1.27%  "
$ws.Range("D51").Value = "2.226.98"
$ws.Range("E51").Value = "  +0.58%  "
